# Add the newest dividend record (XD/Pay date 29/09/2025, Gross Dividend
# 0.009) to the top of the DividendHistory table, just below the header
# row. The existing history rows shift down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DividendHistory")

# Make room for the new record right under the header row.
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = "29/09/2025"
$ws.Cells.Item(2, 2).Value = "29/09/2025"

# "0.009" must stay a text value (matching every other Gross Dividend
# cell in the column), so force text entry with a leading apostrophe
# instead of letting it be parsed as a number, then drop back to the
# Normal style so no stray number-format is left on the cell.
$ws.Cells.Item(2, 3).Value = "'0.009"
$ws.Cells.Item(2, 3).Style = "Normal"
